$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header that used to say "fantasy points" (column E) is relabeled
# "height"; two new columns are appended after it: F = "weight" and
# G = "fantasy points" (the original metric, now living in its own column).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Data rows 2-13: E gets the new "height" value for every player, F gets the
# new "weight" value for every player, and G gets the value that used to live
# in E (0 for everyone except the last row, which keeps 6.1).
for ($r = 2; $r -le 13; $r++) {
    $oldFantasyPoints = $ws.Cells.Item($r, 5).Value2

    $ws.Cells.Item($r, 5).Value = 6.333333333333333
    $ws.Cells.Item($r, 6).Value = 255
    $ws.Cells.Item($r, 7).Value = $oldFantasyPoints
}
